# Replace the old "Generic Asset" folder value with the new Orchestrator
# folder path throughout the workbook.

$wb = $excel.ActiveWorkbook

$newValue = "Downstream/International Products/IFL Automation"

# Settings sheet: B3 held the old value paired with A3 = OrchestratorFolderPath
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B3").Value = $newValue

# Assets sheet: column C (rows 2-7) held the old value as the Asset's folder
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("C2:C7").Value = $newValue

# Reflect the selection left on the Assets sheet (C2:C7)
$wsAssets.Activate()
$wsAssets.Range("C2:C7").Select()
